# Auto-generated COM script applying the commit "Horarios actualizados Linea 141 - 228"
# Updates three worksheets (LP1912, LP1912-215, 6203-6173) with refreshed scrape data:
#  - header rows 2/3 ("Ultima actualizacion" / "Total filas")
#  - several re-sorted rows within tied arrival-time groups
#  - newly scraped rows appended at the end of LP1912 and LP1912-215

$wb = $excel.ActiveWorkbook

# ---- Worksheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Última actualización: 17:53:46"
$ws.Cells.Item(3, 1).Value = "Total filas: 429"
$ws.Cells.Item(45, 1).Value = "06:52:23"
$ws.Cells.Item(45, 2).Value = "07:16"
$ws.Cells.Item(45, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(45, 4).Value = 24
$ws.Cells.Item(45, 5).Value = "LP1912"
$ws.Cells.Item(46, 1).Value = "05:20:00"
$ws.Cells.Item(46, 2).Value = "07:16"
$ws.Cells.Item(46, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(46, 4).Value = 116
$ws.Cells.Item(46, 5).Value = "LP1912"
$ws.Cells.Item(124, 1).Value = "09:38:04"
$ws.Cells.Item(124, 2).Value = "09:41"
$ws.Cells.Item(124, 3).Value = "14_ABASTO"
$ws.Cells.Item(124, 4).Value = 3
$ws.Cells.Item(124, 5).Value = "LP1912"
$ws.Cells.Item(125, 1).Value = "08:39:56"
$ws.Cells.Item(125, 2).Value = "09:41"
$ws.Cells.Item(125, 3).Value = "215C_EL PATO"
$ws.Cells.Item(125, 4).Value = 62
$ws.Cells.Item(125, 5).Value = "LP1912"
$ws.Cells.Item(126, 1).Value = "09:38:04"
$ws.Cells.Item(126, 2).Value = "09:41"
$ws.Cells.Item(126, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(126, 4).Value = 3
$ws.Cells.Item(126, 5).Value = "LP1912"
$ws.Cells.Item(152, 1).Value = "10:57:58"
$ws.Cells.Item(152, 2).Value = "10:57"
$ws.Cells.Item(152, 3).Value = "17_ROMERO"
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 5).Value = "LP1912"
$ws.Cells.Item(154, 1).Value = "10:28:12"
$ws.Cells.Item(154, 2).Value = "10:57"
$ws.Cells.Item(154, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(154, 4).Value = 29
$ws.Cells.Item(154, 5).Value = "LP1912"
$ws.Cells.Item(206, 1).Value = "11:51:05"
$ws.Cells.Item(206, 2).Value = "12:13"
$ws.Cells.Item(206, 3).Value = "10_OLMOS"
$ws.Cells.Item(206, 4).Value = 22
$ws.Cells.Item(206, 5).Value = "LP1912"
$ws.Cells.Item(207, 1).Value = "11:51:05"
$ws.Cells.Item(207, 2).Value = "12:13"
$ws.Cells.Item(207, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(207, 4).Value = 22
$ws.Cells.Item(207, 5).Value = "LP1912"
$ws.Cells.Item(220, 1).Value = "12:16:51"
$ws.Cells.Item(220, 2).Value = "12:34"
$ws.Cells.Item(220, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(220, 4).Value = 18
$ws.Cells.Item(220, 5).Value = "LP1912"
$ws.Cells.Item(221, 1).Value = "12:16:51"
$ws.Cells.Item(221, 2).Value = "12:34"
$ws.Cells.Item(221, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(221, 4).Value = 18
$ws.Cells.Item(221, 5).Value = "LP1912"
$ws.Cells.Item(246, 1).Value = "11:51:05"
$ws.Cells.Item(246, 2).Value = "13:14"
$ws.Cells.Item(246, 3).Value = "215D_EL PATO"
$ws.Cells.Item(246, 4).Value = 83
$ws.Cells.Item(246, 5).Value = "LP1912"
$ws.Cells.Item(247, 1).Value = "12:44:21"
$ws.Cells.Item(247, 2).Value = "13:14"
$ws.Cells.Item(247, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(247, 4).Value = 30
$ws.Cells.Item(247, 5).Value = "LP1912"
$ws.Cells.Item(354, 1).Value = "14:53:55"
$ws.Cells.Item(354, 2).Value = "16:42"
$ws.Cells.Item(354, 3).Value = "225_GOMEZ"
$ws.Cells.Item(354, 4).Value = 109
$ws.Cells.Item(354, 5).Value = "LP1912"
$ws.Cells.Item(355, 1).Value = "14:53:55"
$ws.Cells.Item(355, 2).Value = "16:42"
$ws.Cells.Item(355, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(355, 4).Value = 109
$ws.Cells.Item(355, 5).Value = "LP1912"
$ws.Cells.Item(380, 1).Value = "16:45:22"
$ws.Cells.Item(380, 2).Value = "17:35"
$ws.Cells.Item(380, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(380, 4).Value = 50
$ws.Cells.Item(380, 5).Value = "LP1912"
$ws.Cells.Item(382, 1).Value = "16:14:52"
$ws.Cells.Item(382, 2).Value = "17:35"
$ws.Cells.Item(382, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(382, 4).Value = 81
$ws.Cells.Item(382, 5).Value = "LP1912"
$ws.Cells.Item(396, 1).Value = "17:39:57"
$ws.Cells.Item(396, 2).Value = "17:52"
$ws.Cells.Item(396, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(396, 4).Value = 13
$ws.Cells.Item(396, 5).Value = "LP1912"
$ws.Cells.Item(397, 1).Value = "16:14:52"
$ws.Cells.Item(397, 2).Value = "17:52"
$ws.Cells.Item(397, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(397, 4).Value = 98
$ws.Cells.Item(397, 5).Value = "LP1912"
$ws.Cells.Item(398, 1).Value = "17:53:46"
$ws.Cells.Item(398, 2).Value = "17:53"
$ws.Cells.Item(398, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(398, 4).Value = 0
$ws.Cells.Item(398, 5).Value = "LP1912"
$ws.Cells.Item(399, 1).Value = "17:53:46"
$ws.Cells.Item(399, 2).Value = "17:53"
$ws.Cells.Item(399, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(399, 4).Value = 0
$ws.Cells.Item(399, 5).Value = "LP1912"
$ws.Cells.Item(400, 1).Value = "17:14:54"
$ws.Cells.Item(400, 2).Value = "17:59"
$ws.Cells.Item(400, 3).Value = "10_OLMOS"
$ws.Cells.Item(400, 4).Value = 45
$ws.Cells.Item(400, 5).Value = "LP1912"
$ws.Cells.Item(401, 1).Value = "17:53:46"
$ws.Cells.Item(401, 2).Value = "18:00"
$ws.Cells.Item(401, 3).Value = "10_OLMOS"
$ws.Cells.Item(401, 4).Value = 7
$ws.Cells.Item(401, 5).Value = "LP1912"
$ws.Cells.Item(402, 1).Value = "17:39:57"
$ws.Cells.Item(402, 2).Value = "18:00"
$ws.Cells.Item(402, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(402, 4).Value = 21
$ws.Cells.Item(402, 5).Value = "LP1912"
$ws.Cells.Item(403, 1).Value = "17:39:57"
$ws.Cells.Item(403, 2).Value = "18:03"
$ws.Cells.Item(403, 3).Value = "17_ROMERO"
$ws.Cells.Item(403, 4).Value = 24
$ws.Cells.Item(403, 5).Value = "LP1912"
$ws.Cells.Item(404, 1).Value = "17:14:54"
$ws.Cells.Item(404, 2).Value = "18:04"
$ws.Cells.Item(404, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(404, 4).Value = 50
$ws.Cells.Item(404, 5).Value = "LP1912"
$ws.Cells.Item(405, 1).Value = "16:14:52"
$ws.Cells.Item(405, 2).Value = "18:04"
$ws.Cells.Item(405, 3).Value = "17_ROMERO"
$ws.Cells.Item(405, 4).Value = 110
$ws.Cells.Item(405, 5).Value = "LP1912"
$ws.Cells.Item(406, 1).Value = "16:52:27"
$ws.Cells.Item(406, 2).Value = "18:08"
$ws.Cells.Item(406, 3).Value = "14_ABASTO"
$ws.Cells.Item(406, 4).Value = 76
$ws.Cells.Item(406, 5).Value = "LP1912"
$ws.Cells.Item(407, 1).Value = "17:53:46"
$ws.Cells.Item(407, 2).Value = "18:10"
$ws.Cells.Item(407, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(407, 4).Value = 17
$ws.Cells.Item(407, 5).Value = "LP1912"
$ws.Cells.Item(408, 1).Value = "17:39:57"
$ws.Cells.Item(408, 2).Value = "18:15"
$ws.Cells.Item(408, 3).Value = "15_ABASTO"
$ws.Cells.Item(408, 4).Value = 36
$ws.Cells.Item(408, 5).Value = "LP1912"
$ws.Cells.Item(409, 1).Value = "17:39:57"
$ws.Cells.Item(409, 2).Value = "18:15"
$ws.Cells.Item(409, 3).Value = "10_OLMOS"
$ws.Cells.Item(409, 4).Value = 36
$ws.Cells.Item(409, 5).Value = "LP1912"
$ws.Cells.Item(410, 1).Value = "17:53:46"
$ws.Cells.Item(410, 2).Value = "18:16"
$ws.Cells.Item(410, 3).Value = "15_ABASTO"
$ws.Cells.Item(410, 4).Value = 23
$ws.Cells.Item(410, 5).Value = "LP1912"
$ws.Cells.Item(411, 1).Value = "17:53:46"
$ws.Cells.Item(411, 2).Value = "18:16"
$ws.Cells.Item(411, 3).Value = "10_OLMOS"
$ws.Cells.Item(411, 4).Value = 23
$ws.Cells.Item(411, 5).Value = "LP1912"
$ws.Cells.Item(412, 1).Value = "17:14:54"
$ws.Cells.Item(412, 2).Value = "18:20"
$ws.Cells.Item(412, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(412, 4).Value = 66
$ws.Cells.Item(412, 5).Value = "LP1912"
$ws.Cells.Item(413, 1).Value = "16:32:38"
$ws.Cells.Item(413, 2).Value = "18:21"
$ws.Cells.Item(413, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(413, 4).Value = 109
$ws.Cells.Item(413, 5).Value = "LP1912"
$ws.Cells.Item(414, 1).Value = "17:39:57"
$ws.Cells.Item(414, 2).Value = "18:24"
$ws.Cells.Item(414, 3).Value = "14_ABASTO"
$ws.Cells.Item(414, 4).Value = 45
$ws.Cells.Item(414, 5).Value = "LP1912"
$ws.Cells.Item(415, 1).Value = "16:32:38"
$ws.Cells.Item(415, 2).Value = "18:27"
$ws.Cells.Item(415, 3).Value = "215C_EL PATO"
$ws.Cells.Item(415, 4).Value = 115
$ws.Cells.Item(415, 5).Value = "LP1912"
$ws.Cells.Item(416, 1).Value = "16:45:22"
$ws.Cells.Item(416, 2).Value = "18:28"
$ws.Cells.Item(416, 3).Value = "215C_EL PATO"
$ws.Cells.Item(416, 4).Value = 103
$ws.Cells.Item(416, 5).Value = "LP1912"
$ws.Cells.Item(417, 1).Value = "17:14:54"
$ws.Cells.Item(417, 2).Value = "18:31"
$ws.Cells.Item(417, 3).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(417, 4).Value = 77
$ws.Cells.Item(417, 5).Value = "LP1912"
$ws.Cells.Item(418, 1).Value = "16:45:22"
$ws.Cells.Item(418, 2).Value = "18:32"
$ws.Cells.Item(418, 3).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(418, 4).Value = 107
$ws.Cells.Item(418, 5).Value = "LP1912"
$ws.Cells.Item(419, 1).Value = "17:39:57"
$ws.Cells.Item(419, 2).Value = "18:36"
$ws.Cells.Item(419, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(419, 4).Value = 57
$ws.Cells.Item(419, 5).Value = "LP1912"
$ws.Cells.Item(420, 1).Value = "17:53:46"
$ws.Cells.Item(420, 2).Value = "18:40"
$ws.Cells.Item(420, 3).Value = "15_ABASTO"
$ws.Cells.Item(420, 4).Value = 47
$ws.Cells.Item(420, 5).Value = "LP1912"
$ws.Cells.Item(421, 1).Value = "17:14:54"
$ws.Cells.Item(421, 2).Value = "18:47"
$ws.Cells.Item(421, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(421, 4).Value = 93
$ws.Cells.Item(421, 5).Value = "LP1912"
$ws.Cells.Item(422, 1).Value = "16:52:27"
$ws.Cells.Item(422, 2).Value = "18:48"
$ws.Cells.Item(422, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(422, 4).Value = 116
$ws.Cells.Item(422, 5).Value = "LP1912"
$ws.Cells.Item(423, 1).Value = "17:14:54"
$ws.Cells.Item(423, 2).Value = "18:58"
$ws.Cells.Item(423, 3).Value = "215A_EL PATO"
$ws.Cells.Item(423, 4).Value = 104
$ws.Cells.Item(423, 5).Value = "LP1912"
$ws.Cells.Item(424, 1).Value = "17:14:54"
$ws.Cells.Item(424, 2).Value = "19:04"
$ws.Cells.Item(424, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(424, 4).Value = 110
$ws.Cells.Item(424, 5).Value = "LP1912"
$ws.Cells.Item(425, 1).Value = "17:14:54"
$ws.Cells.Item(425, 2).Value = "19:10"
$ws.Cells.Item(425, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(425, 4).Value = 116
$ws.Cells.Item(425, 5).Value = "LP1912"
$ws.Cells.Item(426, 1).Value = "17:39:57"
$ws.Cells.Item(426, 2).Value = "19:16"
$ws.Cells.Item(426, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(426, 4).Value = 97
$ws.Cells.Item(426, 5).Value = "LP1912"
$ws.Cells.Item(427, 1).Value = "17:39:57"
$ws.Cells.Item(427, 2).Value = "19:20"
$ws.Cells.Item(427, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(427, 4).Value = 101
$ws.Cells.Item(427, 5).Value = "LP1912"
$ws.Cells.Item(428, 1).Value = "17:53:46"
$ws.Cells.Item(428, 2).Value = "19:21"
$ws.Cells.Item(428, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(428, 4).Value = 88
$ws.Cells.Item(428, 5).Value = "LP1912"
$ws.Cells.Item(429, 1).Value = "17:39:57"
$ws.Cells.Item(429, 2).Value = "19:29"
$ws.Cells.Item(429, 3).Value = "225_GOMEZ"
$ws.Cells.Item(429, 4).Value = 110
$ws.Cells.Item(429, 5).Value = "LP1912"
$ws.Cells.Item(430, 1).Value = "17:53:46"
$ws.Cells.Item(430, 2).Value = "19:30"
$ws.Cells.Item(430, 3).Value = "225_GOMEZ"
$ws.Cells.Item(430, 4).Value = 97
$ws.Cells.Item(430, 5).Value = "LP1912"
$ws.Cells.Item(431, 1).Value = "17:53:46"
$ws.Cells.Item(431, 2).Value = "19:39"
$ws.Cells.Item(431, 3).Value = "215C_EL PATO"
$ws.Cells.Item(431, 4).Value = 106
$ws.Cells.Item(431, 5).Value = "LP1912"
$ws.Cells.Item(432, 1).Value = "17:53:46"
$ws.Cells.Item(432, 2).Value = "19:50"
$ws.Cells.Item(432, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(432, 4).Value = 117
$ws.Cells.Item(432, 5).Value = "LP1912"
$ws.Cells.Item(433, 1).Value = "17:53:46"
$ws.Cells.Item(433, 2).Value = "19:50"
$ws.Cells.Item(433, 3).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(433, 4).Value = 117
$ws.Cells.Item(433, 5).Value = "LP1912"
$ws.Cells.Item(434, 1).Value = "17:53:46"
$ws.Cells.Item(434, 2).Value = "19:51"
$ws.Cells.Item(434, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(434, 4).Value = 118
$ws.Cells.Item(434, 5).Value = "LP1912"

# ---- Worksheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Última actualización: 17:53:46"
$ws.Cells.Item(3, 1).Value = "Total filas: 43"
$ws.Cells.Item(48, 1).Value = "17:53:46"
$ws.Cells.Item(48, 2).Value = "19:39"
$ws.Cells.Item(48, 3).Value = "215C_EL PATO"
$ws.Cells.Item(48, 4).Value = 106
$ws.Cells.Item(48, 5).Value = "LP1912"

# ---- Worksheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Última actualización: 17:53:46"

